# gr_Mix of complex & simple sentences changed to another model - gr_Complexity of the verb phrase
#
# Appends 7 new rows (234-240) to Sheet1, duplicating the existing
# "Speaking" / gr_r_covp / "Complexity of the verb phrase" Achievement +
# Suggestion rows (rows 2-4 and 47-50), but for the "Writing" test type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column values shared by every new row.
$topCriteria  = "Grammatical range and accuracy"
$criteriaCode = "gr_r_covp"
$criteria     = "Complexity of the verb phrase"
$status       = "Active"

# --- Achievement rows (234-236) -> copy formatting from row 2 (s="1" on G) ---
$achievementTexts = @(
    "You demonstrate a strong command of various tenses, conditionals, and complex sentence structures",
    "You effectively use passive and active voices to convey different nuances in meaning. ",
    "You employ a range of clauses (relative, noun, and adverbial) accurately and appropriately."
)

$r = 234
foreach ($text in $achievementTexts) {
    $ws.Range("A$r").Value = "Writing"
    $ws.Range("B$r").Value = "Achievement"
    $ws.Range("C$r").Value = $topCriteria
    $ws.Range("D$r").Value = $criteriaCode
    $ws.Range("E$r").Value = $criteria
    $ws.Range("F$r").Value = $status
    $ws.Range("G$r").Value = $text

    [void]$ws.Range("G2").Copy()
    [void]$ws.Range("G$r").PasteSpecial(-4122)

    $r++
}

# --- Suggestion rows (237-240) -> copy formatting from row 47 (no explicit style on G) ---
$suggestionTexts = @(
    "Practice using a wider variety of tenses to improve your grammatical range.",
    "Focus on incorporating more complex sentences into your writing and speaking.",
    "Study and practice different types of clauses to enhance sentence complexity.",
    "Use grammar exercises and resources to reinforce your understanding of advanced structures."
)

foreach ($text in $suggestionTexts) {
    $ws.Range("A$r").Value = "Writing"
    $ws.Range("B$r").Value = "Suggestion"
    $ws.Range("C$r").Value = $topCriteria
    $ws.Range("D$r").Value = $criteriaCode
    $ws.Range("E$r").Value = $criteria
    $ws.Range("F$r").Value = $status
    $ws.Range("G$r").Value = $text

    [void]$ws.Range("G47").Copy()
    [void]$ws.Range("G$r").PasteSpecial(-4122)

    $r++
}

$excel.CutCopyMode = 0

# Keep the frozen-pane view & active selection in sync with the new bottom
# of the sheet, matching how Excel would leave the view after scrolling to
# the newly-entered rows.
[void]$ws.Range("F217").Select()
[void]$ws.Range("E240").Select()
